$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix formatting while old layout/styles are still in place ---
# Clear cell contents only (keeps existing styles: row1=header, rows2-51=data)
$ws.Range("A1:E51").ClearContents()

# Row 2 becomes a repeated header row -> give it the header style (copied from row 1)
$ws.Range("A1:E1").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)

# Rows 52-55 are brand new data rows -> give them the data style (copied from row 3)
$ws.Range("A3:E3").Copy()
$ws.Range("A52:E55").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 2: populate cell values. ---
# Order matters for how new shared strings get appended: the newly-appended
# WAREHOUSE_LIST..ROLE_MANAGE rows (46-55) are written before the Vietnamese
# "Quan ly ..." translations are applied to the *_MANAGE rows, matching how the
# workbook was actually produced.

# 2a. Rows 1-45 (use placeholder/English text for the name column of *_MANAGE rows;
#     will be corrected to Vietnamese in step 2c).
$ws.Cells.Item(1, 1).Value = "code"
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "url"
$ws.Cells.Item(1, 4).Value = "method"
$ws.Cells.Item(1, 5).Value = "parent_cd"

$ws.Cells.Item(2, 1).Value = "code"
$ws.Cells.Item(2, 2).Value = "name"
$ws.Cells.Item(2, 3).Value = "url"
$ws.Cells.Item(2, 4).Value = "method"
$ws.Cells.Item(2, 5).Value = "parent_cd"

$ws.Cells.Item(3, 1).Value = "WAREHOUSE_MANAGE"
$ws.Cells.Item(3, 3).Value = "api/warehouse"
$ws.Cells.Item(3, 4).Value = "GET"

$ws.Cells.Item(4, 1).Value = "WAREHOUSE_DETAIL"
$ws.Cells.Item(4, 2).Value = "Detail Warehouse"
$ws.Cells.Item(4, 3).Value = "api/warehouse/{id}"
$ws.Cells.Item(4, 4).Value = "GET"
$ws.Cells.Item(4, 5).Value = "WAREHOUSE_MANAGE"

$ws.Cells.Item(5, 1).Value = "WAREHOUSE_CREATE"
$ws.Cells.Item(5, 2).Value = "Create Warehouse"
$ws.Cells.Item(5, 3).Value = "api/warehouse/create"
$ws.Cells.Item(5, 4).Value = "POST"
$ws.Cells.Item(5, 5).Value = "WAREHOUSE_MANAGE"

$ws.Cells.Item(6, 1).Value = "WAREHOUSE_UPDATE"
$ws.Cells.Item(6, 2).Value = "Update Warehouse"
$ws.Cells.Item(6, 3).Value = "api/warehouse/update/{id}"
$ws.Cells.Item(6, 4).Value = "PUT"
$ws.Cells.Item(6, 5).Value = "WAREHOUSE_MANAGE"

$ws.Cells.Item(7, 1).Value = "WAREHOUSE_DELETE"
$ws.Cells.Item(7, 2).Value = "Delete Warehouse"
$ws.Cells.Item(7, 3).Value = "api/warehouse/delete/{id}"
$ws.Cells.Item(7, 4).Value = "DELETE"
$ws.Cells.Item(7, 5).Value = "WAREHOUSE_MANAGE"

$ws.Cells.Item(8, 1).Value = "SHIPPER_MANAGE"
$ws.Cells.Item(8, 3).Value = "api/shipper"
$ws.Cells.Item(8, 4).Value = "GET"

$ws.Cells.Item(9, 1).Value = "SHIPPER_DETAIL"
$ws.Cells.Item(9, 2).Value = "Detail Shipper"
$ws.Cells.Item(9, 3).Value = "api/shipper/{id}"
$ws.Cells.Item(9, 4).Value = "GET"
$ws.Cells.Item(9, 5).Value = "SHIPPER_MANAGE"

$ws.Cells.Item(10, 1).Value = "SHIPPER_CREATE"
$ws.Cells.Item(10, 2).Value = "Create Shipper"
$ws.Cells.Item(10, 3).Value = "api/shipper/create"
$ws.Cells.Item(10, 4).Value = "POST"
$ws.Cells.Item(10, 5).Value = "SHIPPER_MANAGE"

$ws.Cells.Item(11, 1).Value = "SHIPPER_UPDATE"
$ws.Cells.Item(11, 2).Value = "Update Shipper"
$ws.Cells.Item(11, 3).Value = "api/shipper/update/{id}"
$ws.Cells.Item(11, 4).Value = "PUT"
$ws.Cells.Item(11, 5).Value = "SHIPPER_MANAGE"

$ws.Cells.Item(12, 1).Value = "SHIPPER_DELETE"
$ws.Cells.Item(12, 2).Value = "Delete Shipper"
$ws.Cells.Item(12, 3).Value = "api/shipper/delete/{id}"
$ws.Cells.Item(12, 4).Value = "DELETE"
$ws.Cells.Item(12, 5).Value = "SHIPPER_MANAGE"

$ws.Cells.Item(13, 1).Value = "CUSTOMER_MANAGE"
$ws.Cells.Item(13, 3).Value = "api/customer"
$ws.Cells.Item(13, 4).Value = "GET"

$ws.Cells.Item(14, 1).Value = "CUSTOMER_DETAIL"
$ws.Cells.Item(14, 2).Value = "Detail Customer"
$ws.Cells.Item(14, 3).Value = "api/customer/{id}"
$ws.Cells.Item(14, 4).Value = "GET"
$ws.Cells.Item(14, 5).Value = "CUSTOMER_MANAGE"

$ws.Cells.Item(15, 1).Value = "CUSTOMER_CREATE"
$ws.Cells.Item(15, 2).Value = "Create Customer"
$ws.Cells.Item(15, 3).Value = "api/customer/create"
$ws.Cells.Item(15, 4).Value = "POST"
$ws.Cells.Item(15, 5).Value = "CUSTOMER_MANAGE"

$ws.Cells.Item(16, 1).Value = "CUSTOMER_UPDATE"
$ws.Cells.Item(16, 2).Value = "Update Customer"
$ws.Cells.Item(16, 3).Value = "api/customer/update/{id}"
$ws.Cells.Item(16, 4).Value = "PUT"
$ws.Cells.Item(16, 5).Value = "CUSTOMER_MANAGE"

$ws.Cells.Item(17, 1).Value = "CUSTOMER_DELETE"
$ws.Cells.Item(17, 2).Value = "Delete Customer"
$ws.Cells.Item(17, 3).Value = "api/customer/delete/{id}"
$ws.Cells.Item(17, 4).Value = "DELETE"
$ws.Cells.Item(17, 5).Value = "CUSTOMER_MANAGE"

$ws.Cells.Item(18, 1).Value = "ORDER_MANAGE"
$ws.Cells.Item(18, 3).Value = "api/order"
$ws.Cells.Item(18, 4).Value = "GET"

$ws.Cells.Item(19, 1).Value = "ORDER_DETAIL"
$ws.Cells.Item(19, 2).Value = "Detail Order"
$ws.Cells.Item(19, 3).Value = "api/order/{id}"
$ws.Cells.Item(19, 4).Value = "GET"
$ws.Cells.Item(19, 5).Value = "ORDER_MANAGE"

$ws.Cells.Item(20, 1).Value = "ORDER_CREATE"
$ws.Cells.Item(20, 2).Value = "Create Order"
$ws.Cells.Item(20, 3).Value = "api/order/create"
$ws.Cells.Item(20, 4).Value = "POST"
$ws.Cells.Item(20, 5).Value = "ORDER_MANAGE"

$ws.Cells.Item(21, 1).Value = "ORDER_UPDATE"
$ws.Cells.Item(21, 2).Value = "Update Order"
$ws.Cells.Item(21, 3).Value = "api/order/update/{id}"
$ws.Cells.Item(21, 4).Value = "PUT"
$ws.Cells.Item(21, 5).Value = "ORDER_MANAGE"

$ws.Cells.Item(22, 1).Value = "ORDER_DELETE"
$ws.Cells.Item(22, 2).Value = "Delete Order"
$ws.Cells.Item(22, 3).Value = "api/order/delete/{id}"
$ws.Cells.Item(22, 4).Value = "DELETE"
$ws.Cells.Item(22, 5).Value = "ORDER_MANAGE"

$ws.Cells.Item(23, 1).Value = "INVOICE_MANAGE"
$ws.Cells.Item(23, 3).Value = "api/invoice"
$ws.Cells.Item(23, 4).Value = "GET"

$ws.Cells.Item(24, 1).Value = "INVOICE_DETAIL"
$ws.Cells.Item(24, 2).Value = "Detail Invoice"
$ws.Cells.Item(24, 3).Value = "api/invoice/{id}"
$ws.Cells.Item(24, 4).Value = "GET"
$ws.Cells.Item(24, 5).Value = "INVOICE_MANAGE"

$ws.Cells.Item(25, 1).Value = "INVOICE_CREATE"
$ws.Cells.Item(25, 2).Value = "Create Invoice"
$ws.Cells.Item(25, 3).Value = "api/invoice/create"
$ws.Cells.Item(25, 4).Value = "POST"
$ws.Cells.Item(25, 5).Value = "INVOICE_MANAGE"

$ws.Cells.Item(26, 1).Value = "INVOICE_UPDATE"
$ws.Cells.Item(26, 2).Value = "Update Invoice"
$ws.Cells.Item(26, 3).Value = "api/invoice/update/{id}"
$ws.Cells.Item(26, 4).Value = "PUT"
$ws.Cells.Item(26, 5).Value = "INVOICE_MANAGE"

$ws.Cells.Item(27, 1).Value = "INVOICE_DELETE"
$ws.Cells.Item(27, 2).Value = "Delete Invoice"
$ws.Cells.Item(27, 3).Value = "api/invoice/delete/{id}"
$ws.Cells.Item(27, 4).Value = "DELETE"
$ws.Cells.Item(27, 5).Value = "INVOICE_MANAGE"

$ws.Cells.Item(28, 1).Value = "INVOICE_EXPORT"
$ws.Cells.Item(28, 2).Value = "Export Invoice"
$ws.Cells.Item(28, 3).Value = "api/invoice/export-invoice"
$ws.Cells.Item(28, 4).Value = "GET"
$ws.Cells.Item(28, 5).Value = "INVOICE_MANAGE"

$ws.Cells.Item(29, 1).Value = "CARTON_MANAGE"
$ws.Cells.Item(29, 3).Value = "api/carton"
$ws.Cells.Item(29, 4).Value = "GET"

$ws.Cells.Item(30, 1).Value = "CARTON_DETAIL"
$ws.Cells.Item(30, 2).Value = "Detail Carton"
$ws.Cells.Item(30, 3).Value = "api/carton/{id}"
$ws.Cells.Item(30, 4).Value = "GET"
$ws.Cells.Item(30, 5).Value = "CARTON_MANAGE"

$ws.Cells.Item(31, 1).Value = "CARTON_CREATE"
$ws.Cells.Item(31, 2).Value = "Create Carton"
$ws.Cells.Item(31, 3).Value = "api/carton/create"
$ws.Cells.Item(31, 4).Value = "POST"
$ws.Cells.Item(31, 5).Value = "CARTON_MANAGE"

$ws.Cells.Item(32, 1).Value = "CARTON_UPDATE"
$ws.Cells.Item(32, 2).Value = "Update Carton"
$ws.Cells.Item(32, 3).Value = "api/carton/update/{id}"
$ws.Cells.Item(32, 4).Value = "PUT"
$ws.Cells.Item(32, 5).Value = "CARTON_MANAGE"

$ws.Cells.Item(33, 1).Value = "CARTON_DELETE"
$ws.Cells.Item(33, 2).Value = "Delete Carton"
$ws.Cells.Item(33, 3).Value = "api/carton/delete/{id}"
$ws.Cells.Item(33, 4).Value = "DELETE"
$ws.Cells.Item(33, 5).Value = "CARTON_MANAGE"

$ws.Cells.Item(34, 1).Value = "PRODUCT_MANAGE"
$ws.Cells.Item(34, 3).Value = "api/product"
$ws.Cells.Item(34, 4).Value = "GET"

$ws.Cells.Item(35, 1).Value = "PRODUCT_DETAIL"
$ws.Cells.Item(35, 2).Value = "Detail Product"
$ws.Cells.Item(35, 3).Value = "api/product/{id}"
$ws.Cells.Item(35, 4).Value = "GET"
$ws.Cells.Item(35, 5).Value = "PRODUCT_MANAGE"

$ws.Cells.Item(36, 1).Value = "PRODUCT_CREATE"
$ws.Cells.Item(36, 2).Value = "Create Product"
$ws.Cells.Item(36, 3).Value = "api/product/create"
$ws.Cells.Item(36, 4).Value = "POST"
$ws.Cells.Item(36, 5).Value = "PRODUCT_MANAGE"

$ws.Cells.Item(37, 1).Value = "PRODUCT_UPDATE"
$ws.Cells.Item(37, 2).Value = "Update Product"
$ws.Cells.Item(37, 3).Value = "api/product/update/{id}"
$ws.Cells.Item(37, 4).Value = "PUT"
$ws.Cells.Item(37, 5).Value = "PRODUCT_MANAGE"

$ws.Cells.Item(38, 1).Value = "PRODUCT_DELETE"
$ws.Cells.Item(38, 2).Value = "Delete Product"
$ws.Cells.Item(38, 3).Value = "api/product/delete/{id}"
$ws.Cells.Item(38, 4).Value = "DELETE"
$ws.Cells.Item(38, 5).Value = "PRODUCT_MANAGE"

$ws.Cells.Item(39, 1).Value = "PRODUCT_EXPORT"
$ws.Cells.Item(39, 2).Value = "Export Product"
$ws.Cells.Item(39, 3).Value = "api/product/export-product"
$ws.Cells.Item(39, 4).Value = "GET"
$ws.Cells.Item(39, 5).Value = "PRODUCT_MANAGE"

$ws.Cells.Item(40, 1).Value = "PRODUCT_IMPORT"
$ws.Cells.Item(40, 2).Value = "Import Product"
$ws.Cells.Item(40, 3).Value = "api/product/import-product"
$ws.Cells.Item(40, 4).Value = "POST"
$ws.Cells.Item(40, 5).Value = "PRODUCT_MANAGE"

$ws.Cells.Item(41, 1).Value = "FUNCTION_MANAGE"
$ws.Cells.Item(41, 3).Value = "api/function"
$ws.Cells.Item(41, 4).Value = "GET"

$ws.Cells.Item(42, 1).Value = "FUNCTION_DETAIL"
$ws.Cells.Item(42, 2).Value = "Detail Function"
$ws.Cells.Item(42, 3).Value = "api/function/{id}"
$ws.Cells.Item(42, 4).Value = "GET"
$ws.Cells.Item(42, 5).Value = "FUNCTION_MANAGE"

$ws.Cells.Item(43, 1).Value = "FUNCTION_CREATE"
$ws.Cells.Item(43, 2).Value = "Create Function"
$ws.Cells.Item(43, 3).Value = "api/function/create"
$ws.Cells.Item(43, 4).Value = "POST"
$ws.Cells.Item(43, 5).Value = "FUNCTION_MANAGE"

$ws.Cells.Item(44, 1).Value = "FUNCTION_UPDATE"
$ws.Cells.Item(44, 2).Value = "Update Function"
$ws.Cells.Item(44, 3).Value = "api/function/update/{id}"
$ws.Cells.Item(44, 4).Value = "PUT"
$ws.Cells.Item(44, 5).Value = "FUNCTION_MANAGE"

$ws.Cells.Item(45, 1).Value = "FUNCTION_DELETE"
$ws.Cells.Item(45, 2).Value = "Delete Function"
$ws.Cells.Item(45, 3).Value = "api/function/delete/{id}"
$ws.Cells.Item(45, 4).Value = "DELETE"
$ws.Cells.Item(45, 5).Value = "FUNCTION_MANAGE"

# 2b. Rows 46-55 (brand-new rows at the bottom of the sheet).
$ws.Cells.Item(46, 1).Value = "WAREHOUSE_LIST"
$ws.Cells.Item(46, 2).Value = "List Warehouse"
$ws.Cells.Item(46, 3).Value = "api/warehouse"
$ws.Cells.Item(46, 4).Value = "GET"
$ws.Cells.Item(46, 5).Value = "WAREHOUSE_MANAGE"

$ws.Cells.Item(47, 1).Value = "SHIPPER_LIST"
$ws.Cells.Item(47, 2).Value = "List Shipper"
$ws.Cells.Item(47, 3).Value = "api/shipper"
$ws.Cells.Item(47, 4).Value = "GET"
$ws.Cells.Item(47, 5).Value = "SHIPPER_MANAGE"

$ws.Cells.Item(48, 1).Value = "CUSTOMER_LIST"
$ws.Cells.Item(48, 2).Value = "List Customer"
$ws.Cells.Item(48, 3).Value = "api/customer"
$ws.Cells.Item(48, 4).Value = "GET"
$ws.Cells.Item(48, 5).Value = "CUSTOMER_MANAGE"

$ws.Cells.Item(49, 1).Value = "ORDER_LIST"
$ws.Cells.Item(49, 2).Value = "List Order"
$ws.Cells.Item(49, 3).Value = "api/order"
$ws.Cells.Item(49, 4).Value = "GET"
$ws.Cells.Item(49, 5).Value = "ORDER_MANAGE"

$ws.Cells.Item(50, 1).Value = "INVOICE_LIST"
$ws.Cells.Item(50, 2).Value = "List Invoice"
$ws.Cells.Item(50, 3).Value = "api/invoice"
$ws.Cells.Item(50, 4).Value = "GET"
$ws.Cells.Item(50, 5).Value = "INVOICE_MANAGE"

$ws.Cells.Item(51, 1).Value = "CARTON_LIST"
$ws.Cells.Item(51, 2).Value = "List Carton"
$ws.Cells.Item(51, 3).Value = "api/carton"
$ws.Cells.Item(51, 4).Value = "GET"
$ws.Cells.Item(51, 5).Value = "CARTON_MANAGE"

$ws.Cells.Item(52, 1).Value = "PRODUCT_LIST"
$ws.Cells.Item(52, 2).Value = "List Product"
$ws.Cells.Item(52, 3).Value = "api/product"
$ws.Cells.Item(52, 4).Value = "GET"
$ws.Cells.Item(52, 5).Value = "PRODUCT_MANAGE"

$ws.Cells.Item(53, 1).Value = "FUNCTION_LIST"
$ws.Cells.Item(53, 2).Value = "List Function"
$ws.Cells.Item(53, 3).Value = "api/function"
$ws.Cells.Item(53, 4).Value = "GET"
$ws.Cells.Item(53, 5).Value = "FUNCTION_MANAGE"

$ws.Cells.Item(54, 1).Value = "USER_LIST"
$ws.Cells.Item(54, 2).Value = "List User"
$ws.Cells.Item(54, 3).Value = "api/user"
$ws.Cells.Item(54, 4).Value = "GET"
$ws.Cells.Item(54, 5).Value = "USER_MANAGE"

$ws.Cells.Item(55, 1).Value = "ROLE_LIST"
$ws.Cells.Item(55, 2).Value = "List Role"
$ws.Cells.Item(55, 3).Value = "api/role"
$ws.Cells.Item(55, 4).Value = "GET"
$ws.Cells.Item(55, 5).Value = "ROLE_MANAGE"

# 2c. Finally, set the Vietnamese names for the *_MANAGE rows.
$ws.Cells.Item(3, 2).Value = "Quản lý kho"
$ws.Cells.Item(8, 2).Value = "Quản lý đơn vị vận chuyển"
$ws.Cells.Item(13, 2).Value = "Quản lý khách hàng"
$ws.Cells.Item(18, 2).Value = "Quản lý đơn hàng"
$ws.Cells.Item(23, 2).Value = "Quản lý hoá đơn"
$ws.Cells.Item(29, 2).Value = "Quản lý thùng"
$ws.Cells.Item(34, 2).Value = "Quản lý sản phẩm"
$ws.Cells.Item(41, 2).Value = "Quản lý phân quyền"

# --- Step 3: restore sheet view (scroll position / selection) ---
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Range("B41").Select()